# Fill in the Day 9 (rows 32-35) and Day 10 (rows 38-41) test-case summary
# numbers that were still blank, and move the sheet's viewport/selection
# down to where the newly entered data lives.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day 9 block: Written / Execution / Review counts
$ws.Range("C33").Value = 7042
$ws.Range("C34").Value = 2402
$ws.Range("C35").Value = 2402

# Day 10 block: Written / Execution / Review counts
$ws.Range("C39").Value = 7045
$ws.Range("C40").Value = 2402
$ws.Range("C41").Value = 2402

# Scroll the window so the edited rows are in view and select the same
# cell Excel left active after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H39").Select()
